$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new "second round of pilot" answers (columns E and F) ---
$ws.Range("E2").Value = "yes, yes"
$ws.Range("F2").Value = "yes, yes"

$ws.Range("E3").Value = "sound quality was disturbung, "
$ws.Range("F3").Value = "not natural piano, the sound was not distrubing"

$ws.Range("E4").Value = "regular, easy to follow, predicting the changes"
$ws.Range("F4").Value = "easy to follow, "

$ws.Range("E5").Value = "difficult, fa"
$ws.Range("F5").Value = "addapted, but hard at the beginning, more on the harder side, "

$ws.Range("E6").Value = "sound more enjoyable"
$ws.Range("F6").Value = "quick learning phase, where the light is"

# Stray note cell further to the right on row 7
$ws.Range("J7").Value = "s"

# New comment cell with its own boxed (left/right border) style, placed at M2
$ws.Range("F2").Copy() | Out-Null
$ws.Range("M2").PasteSpecial(-4122) | Out-Null
$ws.Range("M2").Value = "percpetion bias for different frequencies and loudness"
$ws.Range("M2").Borders.Item(8).LineStyle = 0
$ws.Range("M2").Borders.Item(9).LineStyle = 0

# --- Row 2 grew taller to fit the newly entered text ---
$ws.Rows.Item(2).RowHeight = 102

# --- Move the active selection ---
$ws.Range("B2").Select() | Out-Null
